# Integrated with Test track ALM Client
#
# Applies the workbook changes:
#  - Credentials sheet gains two new UserType rows ("hr", "RM User")
#  - Common / RM_01 sheets get an updated Emp ID value in B2
#  - Selection / active-sheet state is updated to match the new workbook view

$wb = $excel.ActiveWorkbook

$wsCredentials = $wb.Worksheets.Item("Credentials")
$wsCommon      = $wb.Worksheets.Item("Common")
$wsRM01        = $wb.Worksheets.Item("RM_01")

# --- Data edits -----------------------------------------------------------

# New UserType entries on the Credentials sheet. Written A4 then A3 so the
# shared-string table picks up "hr" before "RM User".
$wsCredentials.Range("A4").Value = "hr"
$wsCredentials.Range("A3").Value = "RM User"

# Updated Emp ID values on the other two sheets.
$wsCommon.Range("B2").Value = 22302
$wsRM01.Range("B2").Value = 22302

# --- Selection / active sheet state ----------------------------------------

$wsCredentials.Range("E3").Select()
$wsRM01.Range("B6").Select()

$wsCommon.Activate()
$wsCommon.Range("B2").Select()
